$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Bag method results) appended below existing rows.
$rows = @(
    @{ A = 42602.514305555553; B = "Bag"; C = 5547; D = 4881; E = 289; F = 38; G = 15; H = 71; I = 28; J = 0; K = 1; L = 0; M = 100 },
    @{ A = 42602.516192129631; B = "Bag"; C = 957;  D = 2073; E = 237; F = 33; G = 14; H = 69; I = 29; J = 0; K = 1; L = 0; M = 100 },
    @{ A = 42602.517083333332; B = "Bag"; C = 3531; D = 2454; E = 289; F = 38; G = 15; H = 71; I = 28; J = 0; K = 1; L = 0; M = 100 }
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy the date format (style index) from an existing date cell in
    # column A instead of creating a brand-new number format entry.
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($r, 1).Value = $data.A

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
}
